$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2162.205
$ws.Range("I15").Value = 2162.205
$ws.Range("K15").Value = 6486.615
$ws.Range("M15").Value = -6317.615
$ws.Range("H28").Value = 1038.1666
$ws.Range("I28").Value = 773.8
$ws.Range("K28").Value = 773.8
$ws.Range("M28").Value = -288.8
$ws.Range("H40").Value = 4000
$ws.Range("H43").Value = 3755
$ws.Range("I43").Value = 3749.7778
$ws.Range("J43").Value = 3760.875
$ws.Range("K43").Value = 3749.7778
$ws.Range("L43").Value = 3760.875
$ws.Range("M43").Value = -3680.7778
$ws.Range("N43").Value = -3898.875
$ws.Range("H64").Value = 5480.4287
$ws.Range("I64").Value = 3870
$ws.Range("J64").Value = 5748.8335
$ws.Range("K64").Value = 3870
$ws.Range("L64").Value = 5748.8335
$ws.Range("M64").Value = -3622
$ws.Range("N64").Value = -6244.8335
$ws.Range("H67").Value = 5480.4287
$ws.Range("I67").Value = 3870
$ws.Range("J67").Value = 5748.8335
$ws.Range("K67").Value = 3870
$ws.Range("L67").Value = 5748.8335
$ws.Range("M67").Value = -3012
$ws.Range("N67").Value = -7464.8335
$ws.Range("H74").Value = 10195.586
$ws.Range("I74").Value = 10381.143
$ws.Range("K74").Value = 10381.143
$ws.Range("M74").Value = -9445.143
$ws.Range("H76").Value = 4080.6316
$ws.Range("I76").Value = 4002.75
$ws.Range("J76").Value = 4214.143
$ws.Range("K76").Value = 4002.75
$ws.Range("L76").Value = 4214.143
$ws.Range("M76").Value = -3687.75
$ws.Range("N76").Value = -4844.143
$ws.Range("H77").Value = 10195.586
$ws.Range("I77").Value = 10381.143
$ws.Range("K77").Value = 51905.715
$ws.Range("M77").Value = -47225.715
$ws.Range("H79").Value = 4080.6316
$ws.Range("I79").Value = 4002.75
$ws.Range("J79").Value = 4214.143
$ws.Range("K79").Value = 4002.75
$ws.Range("L79").Value = 4214.143
$ws.Range("M79").Value = -2910.75
$ws.Range("N79").Value = -6398.143
$ws.Range("H137").Value = 5564011.5
$ws.Range("I137").Value = 11113955
$ws.Range("K137").Value = 33341865
$ws.Range("M137").Value = -33339315
$ws.Range("H138").Value = 6909.891
$ws.Range("I138").Value = 6184.1333
$ws.Range("J138").Value = 7261.0645
$ws.Range("K138").Value = 18552.3999
$ws.Range("L138").Value = 21783.1935
$ws.Range("M138").Value = -13412.3999
$ws.Range("N138").Value = -32063.1935

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3344.6667
$ws.Range("I45").Value = 3092.25
$ws.Range("K45").Value = 3092.25
$ws.Range("M45").Value = -2715.25
$ws.Range("H61").Value = 4600010.5
$ws.Range("J61").Value = 10003450
$ws.Range("L61").Value = 10003450
$ws.Range("N61").Value = -10003874
$ws.Range("H122").Value = 1755.875
$ws.Range("I122").Value = 1387
$ws.Range("K122").Value = 4161
$ws.Range("M122").Value = -1711
$ws.Range("H136").Value = 4600010.5
$ws.Range("J136").Value = 10003450
$ws.Range("L136").Value = 30010350
$ws.Range("N136").Value = -30015450

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 241149.05
$ws.Range("J20").Value = 5204.222
$ws.Range("L20").Value = 5204.222
$ws.Range("N20").Value = -5698.222
$ws.Range("H35").Value = 35058
$ws.Range("J35").Value = 35058
$ws.Range("L35").Value = 35058
$ws.Range("N35").Value = -35678
$ws.Range("H134").Value = 3571857.5
$ws.Range("I134").Value = 2860871
$ws.Range("K134").Value = 8582613
$ws.Range("M134").Value = -8580078

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 441984.53
$ws.Range("I31").Value = 880805.9
$ws.Range("J31").Value = 3163.1667
$ws.Range("K31").Value = 880805.9
$ws.Range("L31").Value = 3163.1667
$ws.Range("M31").Value = -880510.9
$ws.Range("N31").Value = -3753.1667
$ws.Range("H34").Value = 441984.53
$ws.Range("I34").Value = 880805.9
$ws.Range("J34").Value = 3163.1667
$ws.Range("K34").Value = 880805.9
$ws.Range("L34").Value = 3163.1667
$ws.Range("M34").Value = -880603.9
$ws.Range("N34").Value = -3567.1667
$ws.Range("H62").Value = 5374
$ws.Range("I62").Value = 4832
$ws.Range("K62").Value = 4832
$ws.Range("M62").Value = -4208
$ws.Range("H65").Value = 5374
$ws.Range("I65").Value = 4832
$ws.Range("K65").Value = 24160
$ws.Range("M65").Value = -21040
$ws.Range("H94").Value = 17900.834
$ws.Range("I94").Value = 33999
$ws.Range("K94").Value = 33999
$ws.Range("M94").Value = -33548
$ws.Range("H105").Value = 6742.222
$ws.Range("I105").Value = 6742.222
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 6742.222
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -4995.222
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 1502.2858
$ws.Range("I107").Value = 1504.75
$ws.Range("K107").Value = 1504.75
$ws.Range("M107").Value = 415.25

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1340781
$ws.Range("J5").Value = 1292647.6
$ws.Range("L5").Value = 3877942.8
$ws.Range("N5").Value = -3878166.8
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H114").Value = 1112.25
$ws.Range("J114").Value = 2999.75
$ws.Range("L114").Value = 8999.25
$ws.Range("N114").Value = -15507.25
$ws.Range("H135").Value = 1340781
$ws.Range("J135").Value = 1292647.6
$ws.Range("L135").Value = 11633828.4
$ws.Range("N135").Value = -11638898.4

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16957.375
$ws.Range("I70").Value = 30701.857
$ws.Range("J70").Value = 6267.222
$ws.Range("K70").Value = 30701.857
$ws.Range("L70").Value = 6267.222
$ws.Range("M70").Value = -30431.857
$ws.Range("N70").Value = -6807.222
$ws.Range("H73").Value = 16957.375
$ws.Range("I73").Value = 30701.857
$ws.Range("J73").Value = 6267.222
$ws.Range("K73").Value = 30701.857
$ws.Range("L73").Value = 6267.222
$ws.Range("M73").Value = -29765.857
$ws.Range("N73").Value = -8139.222
$ws.Range("H98").Value = 29143.5
$ws.Range("J98").Value = 29143.5
$ws.Range("L98").Value = 29143.5
$ws.Range("N98").Value = -35133.5
$ws.Range("H113").Value = 4689.857
$ws.Range("I113").Value = 4643.6
$ws.Range("J113").Value = 4805.5
$ws.Range("K113").Value = 4643.6
$ws.Range("L113").Value = 4805.5
$ws.Range("M113").Value = -2473.6
$ws.Range("N113").Value = -9145.5
$ws.Range("H132").Value = 14608.207
$ws.Range("I132").Value = 8962.125
$ws.Range("J132").Value = 41709.4
$ws.Range("K132").Value = 26886.375
$ws.Range("L132").Value = 125128.2
$ws.Range("M132").Value = -24356.375
$ws.Range("N132").Value = -130188.2

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2553.5386
$ws.Range("I22").Value = 1679.8
$ws.Range("K22").Value = 1679.8
$ws.Range("M22").Value = -1384.8
$ws.Range("H27").Value = 2553.5386
$ws.Range("I27").Value = 1679.8
$ws.Range("K27").Value = 1679.8
$ws.Range("M27").Value = -1572.8
$ws.Range("H40").Value = 4636.4287
$ws.Range("I40").Value = 4492.9165
$ws.Range("J40").Value = 5497.5
$ws.Range("K40").Value = 4492.9165
$ws.Range("L40").Value = 5497.5
$ws.Range("M40").Value = -4356.9165
$ws.Range("N40").Value = -5769.5
$ws.Range("H46").Value = 3184.6553
$ws.Range("I46").Value = 2638.8
$ws.Range("K46").Value = 2638.8
$ws.Range("M46").Value = -2450.8
$ws.Range("H93").Value = 1277.9697
$ws.Range("I93").Value = 875.4761999999999
$ws.Range("J93").Value = 1982.3334
$ws.Range("K93").Value = 875.4761999999999
$ws.Range("L93").Value = 1982.3334
$ws.Range("M93").Value = 372.5238000000001
$ws.Range("N93").Value = -4478.3334

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 97734.63
$ws.Range("I81").Value = 1999.3334
$ws.Range("J81").Value = 133635.38
$ws.Range("K81").Value = 3998.6668
$ws.Range("L81").Value = 267270.76
$ws.Range("M81").Value = -2937.6668
$ws.Range("N81").Value = -269392.76
$ws.Range("H84").Value = 97734.63
$ws.Range("I84").Value = 1999.3334
$ws.Range("J84").Value = 133635.38
$ws.Range("K84").Value = 19993.334
$ws.Range("L84").Value = 1336353.8
$ws.Range("M84").Value = -14689.334
$ws.Range("N84").Value = -1346961.8
$ws.Range("H106").Value = 52500
$ws.Range("J106").Value = 52500
$ws.Range("L106").Value = 52500
$ws.Range("N106").Value = -55024
$ws.Range("H122").Value = 52925.22
$ws.Range("I122").Value = 4468.4707
$ws.Range("K122").Value = 13405.4121
$ws.Range("M122").Value = -10955.4121
$ws.Range("H126").Value = 2741.2727
$ws.Range("I126").Value = 2711.6667
$ws.Range("J126").Value = 2874.5
$ws.Range("K126").Value = 8135.000100000001
$ws.Range("L126").Value = 8623.5
$ws.Range("M126").Value = -5665.000100000001
$ws.Range("N126").Value = -13563.5
$ws.Range("H136").Value = 1768313.4
$ws.Range("I136").Value = 1116001.9
$ws.Range("K136").Value = 3348005.7
$ws.Range("M136").Value = -3345455.7
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
